# Auto-generated Excel COM-interop script to apply numeric updates
# to the Louisoix_Profits workbook tables (currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 808866.4399999999
$ws.Range("I2").Value = 1455015.6
$ws.Range("K2").Value = 1455015.6
$ws.Range("M2").Value = -1454902.6
$ws.Range("H18").Value = 849.6
$ws.Range("I18").Value = 750
$ws.Range("K18").Value = 750
$ws.Range("M18").Value = -466
$ws.Range("H32").Value = 2169.75
$ws.Range("I32").Value = 3333
$ws.Range("J32").Value = 2003.5714
$ws.Range("K32").Value = 3333
$ws.Range("L32").Value = 2003.5714
$ws.Range("M32").Value = -3007
$ws.Range("N32").Value = -2655.5714
$ws.Range("H40").Value = 4753
$ws.Range("I40").Value = 4158.8
$ws.Range("J40").Value = 5495.75
$ws.Range("K40").Value = 4158.8
$ws.Range("L40").Value = 5495.75
$ws.Range("M40").Value = -3983.8
$ws.Range("N40").Value = -5845.75
$ws.Range("H94").Value = 4339.4
$ws.Range("I94").Value = 4339.4
$ws.Range("K94").Value = 4339.4
$ws.Range("M94").Value = -3888.4
$ws.Range("H98").Value = 2714.9092
$ws.Range("I98").Value = 1874
$ws.Range("J98").Value = 6499
$ws.Range("K98").Value = 1874
$ws.Range("L98").Value = 6499
$ws.Range("M98").Value = -376
$ws.Range("N98").Value = -9495
$ws.Range("H112").Value = 2809.6365
$ws.Range("J112").Value = 3369.625
$ws.Range("L112").Value = 10108.875
$ws.Range("N112").Value = -12324.875
$ws.Range("H122").Value = 2714.9092
$ws.Range("I122").Value = 1874
$ws.Range("J122").Value = 6499
$ws.Range("K122").Value = 5622
$ws.Range("L122").Value = 19497
$ws.Range("M122").Value = -3172
$ws.Range("N122").Value = -24397
$ws.Range("H138").Value = 3735.4736
$ws.Range("I138").Value = 3696.25
$ws.Range("J138").Value = 3753.577
$ws.Range("K138").Value = 11088.75
$ws.Range("L138").Value = 11260.731
$ws.Range("M138").Value = -5948.75
$ws.Range("N138").Value = -21540.731
$ws.Range("H141").Value = 2828.0908
$ws.Range("I141").Value = 1684
$ws.Range("K141").Value = 5052
$ws.Range("M141").Value = 128

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22225.666
$ws.Range("I32").Value = 26270.574
$ws.Range("K32").Value = 26270.574
$ws.Range("M32").Value = -25983.574
$ws.Range("H61").Value = 3689.2942
$ws.Range("I61").Value = 3443.8333
$ws.Range("J61").Value = 4278.4
$ws.Range("K61").Value = 3443.8333
$ws.Range("L61").Value = 4278.4
$ws.Range("M61").Value = -3231.8333
$ws.Range("N61").Value = -4702.4
$ws.Range("H97").Value = 6880.4116
$ws.Range("I97").Value = 7511.2856
$ws.Range("K97").Value = 7511.2856
$ws.Range("M97").Value = -7015.2856
$ws.Range("H136").Value = 3689.2942
$ws.Range("I136").Value = 3443.8333
$ws.Range("J136").Value = 4278.4
$ws.Range("K136").Value = 10331.4999
$ws.Range("L136").Value = 12835.2
$ws.Range("M136").Value = -7781.499899999999
$ws.Range("N136").Value = -17935.2
$ws.Range("H139").Value = 80000
$ws.Range("J139").Value = 80000
$ws.Range("L139").Value = 80000
$ws.Range("N139").Value = -90280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 77894
$ws.Range("J81").Value = 77894
$ws.Range("L81").Value = 77894
$ws.Range("N81").Value = -80016
$ws.Range("H84").Value = 77894
$ws.Range("J84").Value = 77894
$ws.Range("L84").Value = 233682
$ws.Range("N84").Value = -244290
$ws.Range("H86").Value = 4022.6667
$ws.Range("I86").Value = 2997.125
$ws.Range("K86").Value = 2997.125
$ws.Range("M86").Value = -1874.125
$ws.Range("H89").Value = 4022.6667
$ws.Range("I89").Value = 2997.125
$ws.Range("K89").Value = 14985.625
$ws.Range("M89").Value = -9369.625
$ws.Range("H99").Value = 132172.38
$ws.Range("I99").Value = 501115
$ws.Range("J99").Value = 9191.5
$ws.Range("K99").Value = 501115
$ws.Range("L99").Value = 9191.5
$ws.Range("M99").Value = -499617
$ws.Range("N99").Value = -12187.5
$ws.Range("H107").Value = 2758.8
$ws.Range("I107").Value = 1664.6666
$ws.Range("K107").Value = 1664.6666
$ws.Range("M107").Value = 255.3334
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 135.71428
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 150
$ws.Range("M4").Value = 12
$ws.Range("N4").Value = -374
$ws.Range("H10").Value = 1224.4546
$ws.Range("I10").Value = 372.5
$ws.Range("J10").Value = 3496.3333
$ws.Range("K10").Value = 372.5
$ws.Range("L10").Value = 3496.3333
$ws.Range("M10").Value = -233.5
$ws.Range("N10").Value = -3774.3333
$ws.Range("H16").Value = 1825.8572
$ws.Range("I16").Value = 1476.5
$ws.Range("J16").Value = 2291.6667
$ws.Range("K16").Value = 1476.5
$ws.Range("L16").Value = 2291.6667
$ws.Range("M16").Value = -1189.5
$ws.Range("N16").Value = -2865.6667
$ws.Range("H58").Value = 77463.64
$ws.Range("I58").Value = 145356.28
$ws.Range("J58").Value = 9571
$ws.Range("K58").Value = 145356.28
$ws.Range("L58").Value = 9571
$ws.Range("M58").Value = -145153.28
$ws.Range("N58").Value = -9977
$ws.Range("H105").Value = 2136.125
$ws.Range("I105").Value = 2136.125
$ws.Range("K105").Value = 2136.125
$ws.Range("M105").Value = -389.125
$ws.Range("H107").Value = 2402.88
$ws.Range("I107").Value = 549.75
$ws.Range("J107").Value = 3274.9412
$ws.Range("K107").Value = 549.75
$ws.Range("L107").Value = 3274.9412
$ws.Range("M107").Value = 1370.25
$ws.Range("N107").Value = -7114.9412
$ws.Range("H113").Value = 1825.8572
$ws.Range("I113").Value = 1476.5
$ws.Range("J113").Value = 2291.6667
$ws.Range("K113").Value = 1476.5
$ws.Range("L113").Value = 2291.6667
$ws.Range("M113").Value = 693.5
$ws.Range("N113").Value = -6631.6667
$ws.Range("H136").Value = 77463.64
$ws.Range("I136").Value = 145356.28
$ws.Range("J136").Value = 9571
$ws.Range("K136").Value = 436068.84
$ws.Range("L136").Value = 28713
$ws.Range("M136").Value = -433518.84
$ws.Range("N136").Value = -33813

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 111.28571
$ws.Range("I2").Value = 49.4
$ws.Range("J2").Value = 266
$ws.Range("K2").Value = 296.4
$ws.Range("L2").Value = 1596
$ws.Range("M2").Value = -183.4
$ws.Range("N2").Value = -1822
$ws.Range("H34").Value = 1244.2
$ws.Range("I34").Value = 510.5
$ws.Range("J34").Value = 1733.3334
$ws.Range("K34").Value = 1531.5
$ws.Range("L34").Value = 5200.0002
$ws.Range("M34").Value = -1447.5
$ws.Range("N34").Value = -5368.0002
$ws.Range("H39").Value = 4403.3477
$ws.Range("I39").Value = 452.5
$ws.Range("J39").Value = 4779.619
$ws.Range("K39").Value = 1357.5
$ws.Range("L39").Value = 14338.857
$ws.Range("M39").Value = -1063.5
$ws.Range("N39").Value = -14926.857
$ws.Range("H55").Value = 6988.778
$ws.Range("J55").Value = 8885
$ws.Range("L55").Value = 26655
$ws.Range("N55").Value = -27009
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 60000
$ws.Range("N94").Value = -61352
$ws.Range("H134").Value = 1273.3334
$ws.Range("I134").Value = 1273.3334
$ws.Range("K134").Value = 3820.0002
$ws.Range("M134").Value = 1249.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 39298.668
$ws.Range("I52").Value = 39948.5
$ws.Range("J52").Value = 37999
$ws.Range("K52").Value = 39948.5
$ws.Range("L52").Value = 37999
$ws.Range("M52").Value = -39689.5
$ws.Range("N52").Value = -38517
$ws.Range("H107").Value = 127157.125
$ws.Range("I107").Value = 167877.83
$ws.Range("J107").Value = 4995
$ws.Range("K107").Value = 167877.83
$ws.Range("L107").Value = 4995
$ws.Range("M107").Value = -165957.83
$ws.Range("N107").Value = -8835

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 5874.875
$ws.Range("J20").Value = 5874.875
$ws.Range("L20").Value = 5874.875
$ws.Range("N20").Value = -6326.875
$ws.Range("H23").Value = 19999
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H68").Value = 3916.5
$ws.Range("I68").Value = 2176.182
$ws.Range("K68").Value = 2176.182
$ws.Range("M68").Value = -1427.182
$ws.Range("H71").Value = 3916.5
$ws.Range("I71").Value = 2176.182
$ws.Range("K71").Value = 10880.91
$ws.Range("M71").Value = -7136.91

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 59999
$ws.Range("J21").Value = 59999
$ws.Range("L21").Value = 59999
$ws.Range("N21").Value = -60469
$ws.Range("H29").Value = 1015999.7
$ws.Range("I29").Value = 1511999.5
$ws.Range("K29").Value = 1511999.5
$ws.Range("M29").Value = -1511709.5
$ws.Range("H35").Value = 59999
$ws.Range("J35").Value = 59999
$ws.Range("L35").Value = 59999
$ws.Range("N35").Value = -60579
$ws.Range("H107").Value = 2335.4614
$ws.Range("I107").Value = 1207
$ws.Range("J107").Value = 4874.5
$ws.Range("K107").Value = 3621
$ws.Range("L107").Value = 14623.5
$ws.Range("M107").Value = -1701
$ws.Range("N107").Value = -18463.5
$ws.Range("H122").Value = 3738.625
$ws.Range("I122").Value = 3738.625
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11215.875
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8765.875
$ws.Range("N122").ClearContents()
